$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J, using the same style as existing headers (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-6 for new columns I and J
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 3

$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 7
